# Updated cryptos list on Sun Nov 24 10:12:00 UTC 2024 with GitHub Actions
#
# All Price (column D) and Volume(1h) (column E) cells on this sheet are
# stored as plain text (they use European "." thousands separators and/or
# padded "  +x.xx%  " strings), so every refresh just overwrites the text
# in place. For the handful of Price cells whose new value happens to look
# like a plain decimal number (e.g. "254.32"), we briefly force a Text
# number format before assigning the value and then clear the format again
# so the cell keeps its original (default) style while still holding the
# text value instead of being auto-converted to a float by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.990.88'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '3.404.70'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '674.74'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  -5.17%  '
$ws.Range("E8").Value = '  -8.72%  '
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '3.402.00'
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +10.51%  '
$ws.Range("D15").Value = '97.732.60'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("E16").Value = '  -3.58%  '
$ws.Range("D17").Value = '4.037.93'
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.81'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +16.43%  '
$ws.Range("D19").Value = '3.403.27'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.579'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +35.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.38'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.07'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.43'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '508.57'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.16%  '
$ws.Range("E25").Value = '  -7.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.55'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '99.72'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '3.587.13'
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("E30").Value = '  +3.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.53'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.194'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +21.93%  '
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.569'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.39'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.44%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +13.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '537.51'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.152'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.63%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.869'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.14%  '
$ws.Range("B45").Value = 'Cosmos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.04'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +14.94%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0434'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("E47").Value = '  -1.40%  '
$ws.Range("E48").Value = '  +14.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +11.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.50'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +8.90%  '
$ws.Range("E51").Value = '  -7.04%  '
